$wb = $excel.ActiveWorkbook

# --- Add "Sheet2" right after the existing "Sheet1" -------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet2"

# --- Header row --------------------------------------------------------------
$ws2.Range("A1").Value = "task"
$ws2.Range("B1").Value = "id"

# --- Task list (A2:A21 = TaskNum1 .. TaskNum20) ------------------------------
for ($i = 1; $i -le 20; $i++) {
    $row = $i + 1
    $ws2.Cells.Item($row, 1).Value = "TaskNum$i"
}

# --- Match the page setup used on the rest of the workbook ------------------
try {
    $ws2.PageSetup.PaperSize = 9
    $ws2.PageSetup.Orientation = 1
} catch {
}

# --- Recreate the view state captured in the saved workbook -----------------
$ws2.Activate()
$ws2.Range("B21").Select()
try {
    $excel.ActiveWindow.ScrollRow = 3
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
}
